# Update Machine Process Condition table (columns L:P) on Sheet1.
# The "Jog Feeding/Bending/Turning Request" rows are each split into a
# "Positif" and "Negatif" variant. This shifts the following
# "Operate Feeding/Bending/Turning Request" rows down by two rows
# (rows 22-24 -> rows 25-27), while "Origin Request" (which already had a
# gap row before it) ends up at row 29 instead of row 26. The Modbus / PLC
# addresses (columns M/N) are renumbered accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - Mode
$ws.Range("L8").Value = "Mode"
$ws.Range("O8").Value = "Manual"
$ws.Range("P8").Value = "Auto"

# Row 9 - Run
$ws.Range("O9").Value = "Stop"

# Row 10 - Alarm
$ws.Range("O10").Value = "No Alarm"

# Row 11 - Reset
$ws.Range("L11").Value = "Reset"
$ws.Range("N11").Value = "M3"
$ws.Range("O11").Value = "No"
$ws.Range("P11").Value = "Request"

# Row 13 - Pressure Die
$ws.Range("L13").Value = "Pressure Die"
$ws.Range("O13").Value = "Open"
$ws.Range("P13").Value = "Close"

# Row 14 - Clamp Die
$ws.Range("L14").Value = "Clamp Die"
$ws.Range("O14").Value = "Open"
$ws.Range("P14").Value = "Close"

# Row 15 - Table Up
$ws.Range("O15").Value = "At Bottom"
$ws.Range("P15").Value = "At Top"

# Row 16 - Table Shift
$ws.Range("O16").Value = "At Left"
$ws.Range("P16").Value = "At Right"

# Row 18 - Jog Enable
$ws.Range("O18").Value = "Disabled"
$ws.Range("P18").Value = "Enabled"

# Row 19 - was "Jog Feeding Request" -> "Jog Feeding Request Positif"
$ws.Range("L19").Value = "Jog Feeding Request Positif"
$ws.Range("O19").Value = "No"
$ws.Range("P19").Value = "Yes"

# Row 20 - was "Jog Bending Request" -> "Jog Feeding Request Negatif"
$ws.Range("L20").Value = "Jog Feeding Request Negatif"
$ws.Range("O20").Value = "No"
$ws.Range("P20").Value = "Yes"

# Row 21 - was "Jog Turning Request" -> "Jog Bending Request Positif"
$ws.Range("L21").Value = "Jog Bending Request Positif"
$ws.Range("O21").Value = "No"
$ws.Range("P21").Value = "Yes"

# Row 22 - was "Operate Feeding Request"/M24 -> "Jog Bending Request Negatif"
$ws.Range("L22").Value = "Jog Bending Request Negatif"
$ws.Range("N22").Value = "M24"
$ws.Range("O22").Value = "No"
$ws.Range("P22").Value = "Yes"

# Row 23 - was "Operate Bending Request"/M25 -> "Jog Turning Request Positif"
$ws.Range("L23").Value = "Jog Turning Request Positif"
$ws.Range("N23").Value = "M25"
$ws.Range("O23").Value = "No"
$ws.Range("P23").Value = "Yes"

# Row 24 - was "Operate Turning Request"/M26 -> "Jog Turning Request Negatif"
$ws.Range("L24").Value = "Jog Turning Request Negatif"
$ws.Range("N24").Value = "M26"
$ws.Range("O24").Value = "No"
$ws.Range("P24").Value = "Yes"

# Row 25 (new) - "Operate Feeding Request"
$ws.Range("L25").Value = "Operate Feeding Request"
$ws.Range("M25").Value = 3099
$ws.Range("N25").Value = "M27"
$ws.Range("O25").Value = "No"
$ws.Range("P25").Value = "Yes"

# Row 26 - was "Origin Request"/3102/M30 -> "Operate Bending Request"
$ws.Range("L26").Value = "Operate Bending Request"
$ws.Range("M26").Value = 3100
$ws.Range("N26").Value = "M28"
$ws.Range("O26").Value = "No"
$ws.Range("P26").Value = "Yes"

# Row 27 (new) - "Operate Turning Request"
$ws.Range("L27").Value = "Operate Turning Request"
$ws.Range("M27").Value = 3101
$ws.Range("N27").Value = "M29"
$ws.Range("O27").Value = "No"
$ws.Range("P27").Value = "Yes"

# Row 29 (new position) - "Origin Request" moves back down here
$ws.Range("L29").Value = "Origin Request"
$ws.Range("M29").Value = 3102
$ws.Range("N29").Value = "M30"
$ws.Range("O29").Value = "No"
$ws.Range("P29").Value = "Yes"

# Update the view selection to roughly match the author's last position.
$ws.Range("L25").Select()
